# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the "Feria Lagunitas de Puerto
# Montt - Mango" data block (row 195), pushing the existing rows 195-218
# down to 196-219, and populate the new row with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 195:218 down to 196:219, keeping their data/formatting intact.
$ws.Rows("195:195").Insert()

# Populate the newly inserted row 195 with the new weekly record.
$ws.Cells.Item(195, 1).Value = 4
$ws.Cells.Item(195, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(195, 3).Value = "Los Lagos"
$ws.Cells.Item(195, 4).Value = 44769
$ws.Cells.Item(195, 4).NumberFormat = $ws.Cells.Item(196, 4).NumberFormat
$ws.Cells.Item(195, 5).Value = 10
$ws.Cells.Item(195, 6).Value = "Fruta"
$ws.Cells.Item(195, 7).Value = 100108
$ws.Cells.Item(195, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(195, 9).Value = 100108002
$ws.Cells.Item(195, 10).Value = "Mango"
$ws.Cells.Item(195, 11).Value = "Sin especificar"
$ws.Cells.Item(195, 12).Value = "Primera"
$ws.Cells.Item(195, 13).Value = 30
$ws.Cells.Item(195, 14).Value = 13000
$ws.Cells.Item(195, 15).Value = 14000
$ws.Cells.Item(195, 16).Value = 13500
$ws.Cells.Item(195, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(195, 18).Value = "Brasil"
$ws.Cells.Item(195, 19).Value = 3375
$ws.Cells.Item(195, 20).Value = 4
